# "1st changes of mifos to finflux"
# Insert a new (blank) column into the "Repayment schedule" sheet right
# before column N (the existing "Late" column), pushing Late / heading /
# Outstanding one column to the right, then make that sheet the active /
# selected tab with the selection parked on R6.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new column before column N (14th column). This shifts N->O,
# O->P, P->Q and keeps each cell's style, pushing the "Late"/"heading"/
# "Outstanding" headers and their data one column over.
$ws.Columns.Item(14).Insert()

# The newly inserted column has no explicit width yet - Excel normally
# carries the width of the column immediately to its left (M) into the
# freshly inserted one.
$ws.Columns.Item(14).ColumnWidth = $ws.Columns.Item(13).ColumnWidth

# Make "Repayment schedule" the active sheet/tab (was "NewLoanInput").
$ws.Activate()
$ws.Range("R6").Select() | Out-Null
